$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so values like
# '72.407.14' are not reinterpreted as numbers/dates by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '72.407.14'
$ws.Range('D3').Value = '4.048.62'
$ws.Range('E3').Value = '  +4.04%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '522.28'
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').Value = '148.61'
$ws.Range('E6').Value = '  +2.92%  '
$ws.Range('D7').Value = '0.718'
$ws.Range('E7').Value = '  +17.27%  '
$ws.Range('D8').Value = '4.039.85'
$ws.Range('E8').Value = '  +4.11%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '0.783'
$ws.Range('E10').Value = '  +8.97%  '
$ws.Range('D11').Value = '0.179'
$ws.Range('E11').Value = '  +4.64%  '
$ws.Range('D12').Value = '0.0000334'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').Value = '48.44'
$ws.Range('E13').Value = '  +15.46%  '
$ws.Range('D14').Value = '11.19'
$ws.Range('E14').Value = '  +9.56%  '
$ws.Range('D15').Value = '4.694.65'
$ws.Range('E15').Value = '  +4.06%  '
$ws.Range('D16').Value = '4.023.10'
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = '21.34'
$ws.Range('E17').Value = '  +6.61%  '
$ws.Range('E18').Value = '  +2.44%  '
$ws.Range('D19').Value = '1.23'
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('D20').Value = '0.134'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').Value = '72.330.50'
$ws.Range('E21').Value = '  +4.56%  '
$ws.Range('D22').Value = '446.73'
$ws.Range('E22').Value = '  +5.53%  '
$ws.Range('D23').Value = '105.10'
$ws.Range('E23').Value = '  +20.06%  '
$ws.Range('E24').Value = '  +6.43%  '
$ws.Range('D25').Value = '15.23'
$ws.Range('E25').Value = '  +7.53%  '
$ws.Range('D26').Value = '4.05'
$ws.Range('E26').Value = '  +0.82%  '
$ws.Range('D27').Value = '11.52'
$ws.Range('E27').Value = '  +1.38%  '
$ws.Range('D28').Value = '11.22'
$ws.Range('E28').Value = '  +6.46%  '
$ws.Range('D29').Value = '38.16'
$ws.Range('E29').Value = '  +4.99%  '
$ws.Range('D30').Value = '5.83'
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('D31').Value = '3.29'
$ws.Range('E31').Value = '  +15.45%  '
$ws.Range('D32').Value = '13.80'
$ws.Range('E32').Value = '  +4.63%  '
$ws.Range('E33').Value = '  +4.11%  '
$ws.Range('D34').Value = '683.07'
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('D35').Value = '6.75'
$ws.Range('E35').Value = '  +14.21%  '
$ws.Range('D36').Value = '67.97'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').Value = '42.72'
$ws.Range('E37').Value = '  +6.79%  '
$ws.Range('E38').Value = '  +2.64%  '
$ws.Range('D39').Value = '0.433'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range('D41').Value = '3.51'
$ws.Range('E41').Value = '  +6.81%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('E43').Value = '  +4.53%  '
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('E45').Value = '  -1.38%  '
$ws.Range('D46').Value = '0.158'
$ws.Range('E46').Value = '  +12.56%  '
$ws.Range('D47').Value = '2.73'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('D48').Value = '3.48'
$ws.Range('E48').Value = '  +2.35%  '
$ws.Range('D49').Value = '9.68'
$ws.Range('E49').Value = '  +13.27%  '
$ws.Range('E50').Value = '  +3.00%  '
$ws.Range('D51').Value = '0.000280'
$ws.Range('E51').Value = '  +3.79%  '
